$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").Value = "'" + '26.238.76'
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = '  -0.70%  '

$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").Value = "'" + '1.657.37'
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = '  -0.96%  '

$ws.Range("E4").Value = '  -0.66%  '

$ws.Range("E5").Value = '  -0.83%  '

$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'" + '0.5236'
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  -2.33%  '

$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").Value = "'" + '1.004'
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = '  -0.62%  '

$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").Value = "'" + '0.2655'
$ws.Range("D8").Style = $style_D8
$ws.Range("E8").Value = '  -0.71%  '

$ws.Range("E9").Value = '  -1.35%  '

$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'" + '20.69'
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  -1.79%  '

$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").Value = "'" + '0.07777'
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = '  -1.04%  '

$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").Value = "'" + '4.560'
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = '  -0.35%  '

$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'" + '1.695.65'
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = '  +1.16%  '

$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'" + '1.884.97'
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = '  -0.94%  '

$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").Value = "'" + '0.5653'
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = '  -0.12%  '

$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").Value = "'" + '0.0₅8105'
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = '  -1.29%  '

$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").Value = "'" + '65.50'
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = '  -1.32%  '

$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'" + '26.228.06'
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  -0.91%  '

$ws.Range("E19").Value = '  -0.63%  '

$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").Value = "'" + '4.729'
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = '  +0.26%  '

$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'" + '193.38'
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = '  -1.98%  '

$ws.Range("E22").Value = '  -0.71%  '

$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'" + '6.033'
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  -0.74%  '

$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'" + '143.90'
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = '  -1.72%  '

$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'" + '0.1203'
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = '  -2.52%  '

$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").Value = "'" + '7.267'
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = '  -0.09%  '

$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").Value = "'" + '16.02'
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = '  -1.33%  '

$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'" + '1.504'
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = '  -0.41%  '

$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").Value = "'" + '0.05605'
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = '  -5.04%  '

$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").Value = "'" + '1.277'
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = '  -1.24%  '

$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'" + '3.504'
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  -2.35%  '

$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").Value = "'" + '3.383'
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = '  +2.03%  '

$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").Value = "'" + '1.590'
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = '  -2.35%  '

$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").Value = "'" + '2.803'
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = '  -1.78%  '

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'" + '0.9443'
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  -3.01%  '

$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").Value = "'" + '2.405'
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = '  -1.02%  '

$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'" + '0.5757'
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("E39").Value = '  -0.84%  '

$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'" + '5.926'
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = '  +0.16%  '

$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'" + '2.587'
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  -0.42%  '

$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'" + '0.8461'
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = '  -2.59%  '

$ws.Range("E43").Value = '  -0.68%  '

$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").Value = "'" + '1.037.89'
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = '  -3.99%  '

$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").Value = "'" + '102.48'
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = '  -1.96%  '

$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").Value = "'" + '1.796.27'
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = '  -0.87%  '

$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'" + '58.44'
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("E48").Value = '  -1.50%  '

$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'" + '1.002'
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'" + '0.05323'
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = '  +2.99%  '

$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").Value = "'" + '0.4354'
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = '  -1.12%  '
